# Apply Betfair Back/Lay odds updates for 2025-11-11
# (values per the authoritative diff of the commit "Atualizando o arquivo XLSX")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.66
$ws.Range("Q2").Value = 1.84
$ws.Range("Y2").Value = 20
$ws.Range("AD2").Value = 24
$ws.Range("I3").Value = 7.4
$ws.Range("L3").Value = 1.29
$ws.Range("X3").Value = 22
$ws.Range("Y3").Value = 28
$ws.Range("Z3").Value = 60
$ws.Range("AA3").Value = 220
$ws.Range("AB3").Value = 11
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 29
$ws.Range("AF3").Value = 12
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 100
$ws.Range("AJ3").Value = 18
$ws.Range("AK3").Value = 20
$ws.Range("AL3").Value = 42
$ws.Range("AN3").Value = 10
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2.08
$ws.Range("H4").Value = 3.5
$ws.Range("J4").Value = 3.95
$ws.Range("Q4").Value = 1.73
$ws.Range("T4").Value = 1.64
$ws.Range("W4").Value = 1.92
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 21
$ws.Range("AB4").Value = 14
$ws.Range("AC4").Value = 11
$ws.Range("AF4").Value = 17.5
$ws.Range("F5").Value = 2.58
$ws.Range("G5").Value = 3.1
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3.65
$ws.Range("K5").Value = 3.25
$ws.Range("L5").Value = 1.44
$ws.Range("N5").Value = 2.5
$ws.Range("O5").Value = 1.46
$ws.Range("P5").Value = 1.5
$ws.Range("Q5").Value = 2.38
$ws.Range("R5").Value = 1.18
$ws.Range("S5").Value = 4.2
$ws.Range("T5").Value = 1.92
$ws.Range("U5").Value = 1.84
$ws.Range("V5").Value = 1.38
$ws.Range("W5").Value = 1.51
$ws.Range("AB5").Value = 11
$ws.Range("AC5").Value = 8.4
$ws.Range("F6").Value = 2.6
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 2.94
$ws.Range("I6").Value = 3.65
$ws.Range("K6").Value = 3.3
$ws.Range("M6").Value = 1.12
$ws.Range("N6").Value = 2.52
$ws.Range("O6").Value = 1.52
$ws.Range("P6").Value = 1.51
$ws.Range("S6").Value = 4.7
$ws.Range("T6").Value = 2.04
$ws.Range("U6").Value = 1.76
$ws.Range("V6").Value = 1.38
$ws.Range("W6").Value = 1.51
$ws.Range("AF6").Value = 21
$ws.Range("G7").Value = 4.9
